$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells: "<name>_old" -> "<name>_FV2410", "<name>_new" -> "<name>_FV2504" ---
$headers = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($headers[$i])_FV2410"
}

# column 11 is "diff" (unchanged)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($headers[$i])_FV2504"
}

# --- Turn the data range into a native Excel Table (ListObject) ---
$range = $ws.Range("A1:U57")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"

# --- Freeze the header row (pane split at row 1) ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
